# Update the "Training Dashboard" sheet: for each training row (3-15),
# decrement the "PERIOD TO EXPIRE" (column H) by 1 day and refresh the
# "LAST UPDATE" (column I) date text from 03-Nov-2025 to 04-Nov-2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 15; $row++) {
    # Column H: PERIOD TO EXPIRE - numeric, decrease by one.
    $hCell = $ws.Cells.Item($row, 8)
    $hCell.Value2 = $hCell.Value2 - 1

    # Column I: LAST UPDATE - stored as literal text (e.g. "04-Nov-2025"),
    # not a real date. Writing the date-like string straight into .Value
    # would make Excel auto-convert it into a date serial number, which
    # would also change the cell's number format. To keep it as plain
    # text (same representation/style as before), build it as a formula
    # that evaluates to the text, then convert that formula result back
    # into a static value in place.
    $iCell = $ws.Cells.Item($row, 9)
    $iCell.Value = '="04-Nov-2025"'
    $iCell.Copy() | Out-Null
    $iCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
